$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.093.92'
$ws.Range('E2').Value = '  -1.02%  '
$ws.Range('D3').Value = '2.091.51'
$ws.Range('E3').Value = '  +3.13%  '
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('E5').Value = '  -1.00%  '
$ws.Range('E6').Value = '  -6.16%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range("D8").Value = "'50.95"
$ws.Range('E8').Value = '  +9.14%  '
$ws.Range("D9").Value = "'60.39"
$ws.Range('E9').Value = '  +4.21%  '
$ws.Range("D10").Value = "'0.369"
$ws.Range('E10').Value = '  -2.92%  '
$ws.Range("D11").Value = "'0.0739"
$ws.Range('E11').Value = '  -3.47%  '
$ws.Range('E12').Value = '  +4.76%  '
$ws.Range("D13").Value = "'15.29"
$ws.Range('E13').Value = '  -1.08%  '
$ws.Range('D14').Value = '2.358.16'
$ws.Range('E14').Value = '  +1.77%  '
$ws.Range("D15").Value = "'0.828"
$ws.Range('E15').Value = '  -1.04%  '
$ws.Range('D16').Value = '2.088.76'
$ws.Range('E16').Value = '  +2.88%  '
$ws.Range("D17").Value = "'5.07"
$ws.Range('E17').Value = '  -2.55%  '
$ws.Range('D18').Value = '37.034.70'
$ws.Range('E18').Value = '  -0.70%  '
$ws.Range("D19").Value = "'72.04"
$ws.Range('E19').Value = '  -4.26%  '
$ws.Range('D20').Value = '0.0₃0821'
$ws.Range('E20').Value = '  -4.54%  '
$ws.Range("D21").Value = "'13.29"
$ws.Range('E21').Value = '  -3.45%  '
$ws.Range("D22").Value = "'239.19"
$ws.Range('E22').Value = '  -5.97%  '
$ws.Range("D23").Value = "'5.21"
$ws.Range('E23').Value = '  -1.31%  '
$ws.Range('E24').Value = '  +0.42%  '
$ws.Range('E25').Value = '  -3.89%  '
$ws.Range("D26").Value = "'169.38"
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range("D27").Value = "'9.14"
$ws.Range('E27').Value = '  +2.37%  '
$ws.Range("D28").Value = "'20.68"
$ws.Range('E28').Value = '  +2.19%  '
$ws.Range("D29").Value = "'2.00"
$ws.Range('E29').Value = '  -5.71%  '
$ws.Range('E30').Value = '  -5.67%  '
$ws.Range("D31").Value = "'1.05"
$ws.Range('E31').Value = '  +17.61%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = "'4.49"
$ws.Range('E32').Value = '  -2.19%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = "'0.0605"
$ws.Range('E33').Value = '  -1.80%  '
$ws.Range('B34').Value = 'Gas'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range("D34").Value = "'20.93"
$ws.Range('E34').Value = '  -6.84%  '
$ws.Range("D35").Value = "'0.0911"
$ws.Range('E35').Value = '  +1.32%  '
$ws.Range('E36').Value = '  +0.00%  '
$ws.Range("D37").Value = "'2.31"
$ws.Range('E37').Value = '  +2.56%  '
$ws.Range("D38").Value = "'4.09"
$ws.Range('E38').Value = '  -5.88%  '
$ws.Range("D39").Value = "'1.82"
$ws.Range('E39').Value = '  -3.67%  '
$ws.Range('E40').Value = '  -9.18%  '
$ws.Range("D41").Value = "'17.80"
$ws.Range('E41').Value = '  +2.01%  '
$ws.Range("D42").Value = "'0.0224"
$ws.Range('E42').Value = '  -2.10%  '
$ws.Range('E43').Value = '  +2.27%  '
$ws.Range("D44").Value = "'98.14"
$ws.Range('E44').Value = '  -5.15%  '
$ws.Range("D45").Value = "'2.77"
$ws.Range('E45').Value = '  -3.88%  '
$ws.Range("D46").Value = "'0.0881"
$ws.Range('E46').Value = '  +3.64%  '
$ws.Range('E47').Value = '  +4.72%  '
$ws.Range('D48').Value = '1.311.09'
$ws.Range('E48').Value = '  -4.44%  '
$ws.Range("D49").Value = "'6.91"
$ws.Range('E49').Value = '  +4.62%  '
$ws.Range('D50').Value = '2.276.02'
$ws.Range('E50').Value = '  +4.41%  '
$ws.Range("D51").Value = "'2.28"
$ws.Range('E51').Value = '  -5.00%  '
